# "New Admin User, Search, Reset"
# Rename Sheet2 -> AdminUserPage, make it the active tab, and populate
# the new admin user row (name / hyperlinked email / role).

$wb = $excel.ActiveWorkbook

# --- rename the second sheet and make it active -------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Name = "AdminUserPage"

# --- new admin user row ---------------------------------------------------
$ws.Range("A1").Value = "John1978"
$ws.Range("B1").Value = "jo@8791"
$ws.Range("C1").Value = "Admin"

# hyperlink the e-mail cell (applies the built-in Hyperlink style too)
$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:jo@8791") | Out-Null

# --- selection / active sheet state matching the authored workbook -------
$ws.Range("C2").Select() | Out-Null
$ws.Activate()
